# Fix the three misspelled/miscapitalized state names so that column A of
# Sheet1 is back in correct alphabetical order (Idaho, Illinois, Indiana,
# Iowa). "Indiana" (row 16) was already spelled correctly, so it is left
# untouched; only the "ldaho", "illinois" and "lowa" typos are corrected.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = "Idaho"
$ws.Range("A15").Value = "Illinois"
$ws.Range("A17").Value = "Iowa"

# Update the view state: scroll the sheet down a bit and move the
# selection to A18, matching the author's new cursor position after
# reviewing/fixing the state list.
$ws.Range("A18").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
